# Auto-generated edit script applying the diff to Sargatanas_Profits workbook
# Updates H/I/J/K/L/M/N profit columns across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 1408.619
$ws.Range("J17").Value2 = 1408.619
$ws.Range("L17").Value2 = 4225.857
$ws.Range("N17").Value2 = -4561.857

$ws.Range("H43").Value2 = 458349.78
$ws.Range("I43").Value2 = 3291.5
$ws.Range("J43").Value2 = 1368466.4
$ws.Range("K43").Value2 = 3291.5
$ws.Range("L43").Value2 = 1368466.4
$ws.Range("M43").Value2 = -3222.5
$ws.Range("N43").Value2 = -1368604.4

$ws.Range("H74").Value2 = 107152280
$ws.Range("I74").Value2 = 214288960
$ws.Range("J74").Value2 = 15607
$ws.Range("K74").Value2 = 214288960
$ws.Range("L74").Value2 = 15607
$ws.Range("M74").Value2 = -214288024
$ws.Range("N74").Value2 = -17479

$ws.Range("H77").Value2 = 107152280
$ws.Range("I77").Value2 = 214288960
$ws.Range("J77").Value2 = 15607
$ws.Range("K77").Value2 = 1071444800
$ws.Range("L77").Value2 = 78035
$ws.Range("M77").Value2 = -1071440120
$ws.Range("N77").Value2 = -87395

$ws.Range("H98").Value2 = 27780258
$ws.Range("I98").Value2 = 28573936
$ws.Range("J98").Value2 = 1500
$ws.Range("K98").Value2 = 28573936
$ws.Range("L98").Value2 = 1500
$ws.Range("M98").Value2 = -28572438
$ws.Range("N98").Value2 = -4496

$ws.Range("H116").Value2 = 11912942
$ws.Range("I116").Value2 = 19236552
$ws.Range("K116").Value2 = 19236552
$ws.Range("M116").Value2 = -19233110

$ws.Range("H122").Value2 = 27780258
$ws.Range("I122").Value2 = 28573936
$ws.Range("J122").Value2 = 1500
$ws.Range("K122").Value2 = 85721808
$ws.Range("L122").Value2 = 4500
$ws.Range("M122").Value2 = -85719358
$ws.Range("N122").Value2 = -9400

$ws.Range("H137").Value2 = 4930.276
$ws.Range("I137").Value2 = 3951.5454
$ws.Range("J137").Value2 = 5528.3887
$ws.Range("K137").Value2 = 11854.6362
$ws.Range("L137").Value2 = 16585.1661
$ws.Range("M137").Value2 = -9304.636200000001
$ws.Range("N137").Value2 = -21685.1661

$ws.Range("H138").Value2 = 2637586.2
$ws.Range("I138").Value2 = 3342.3333
$ws.Range("K138").Value2 = 10026.9999
$ws.Range("M138").Value2 = -4886.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value2 = 14622.556
$ws.Range("I31").Value2 = 3965.8572
$ws.Range("K31").Value2 = 3965.8572
$ws.Range("M31").Value2 = -3671.8572

$ws.Range("H32").Value2 = 1331837
$ws.Range("I32").Value2 = 1331837
$ws.Range("K32").Value2 = 1331837
$ws.Range("M32").Value2 = -1331550

$ws.Range("H45").Value2 = 2411.7778
$ws.Range("I45").Value2 = 1299.5
$ws.Range("K45").Value2 = 1299.5
$ws.Range("M45").Value2 = -922.5

$ws.Range("H74").Value2 = 66630.36
$ws.Range("I74").Value2 = 95279.12
$ws.Range("K74").Value2 = 95279.12
$ws.Range("M74").Value2 = -94405.12

$ws.Range("H77").Value2 = 66630.36
$ws.Range("I77").Value2 = 95279.12
$ws.Range("K77").Value2 = 476395.6
$ws.Range("M77").Value2 = -472027.6

$ws.Range("H82").Value2 = 23019.111
$ws.Range("J82").Value2 = 23019.111
$ws.Range("L82").Value2 = 23019.111
$ws.Range("N82").Value2 = -23741.111

$ws.Range("H85").Value2 = 23019.111
$ws.Range("J85").Value2 = 23019.111
$ws.Range("L85").Value2 = 23019.111
$ws.Range("N85").Value2 = -25515.111

$ws.Range("H102").Value2 = 1617.5483
$ws.Range("I102").Value2 = 1408.7
$ws.Range("K102").Value2 = 1408.7
$ws.Range("M102").Value2 = 213.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value2 = 2248.0417
$ws.Range("I94").Value2 = 1804.5385
$ws.Range("J94").Value2 = 2772.182
$ws.Range("K94").Value2 = 1804.5385
$ws.Range("L94").Value2 = 2772.182
$ws.Range("M94").Value2 = -1353.5385
$ws.Range("N94").Value2 = -3674.182

$ws.Range("H107").Value2 = 59267324
$ws.Range("I107").Value2 = 70378770
$ws.Range("J107").Value2 = 6284
$ws.Range("K107").Value2 = 70378770
$ws.Range("L107").Value2 = 6284
$ws.Range("M107").Value2 = -70376850
$ws.Range("N107").Value2 = -10124

$ws.Range("H113").Value2 = 5495.143
$ws.Range("I113").Value2 = 5495.143
$ws.Range("K113").Value2 = 5495.143
$ws.Range("M113").Value2 = -3325.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 8250.117
$ws.Range("I31").Value2 = 2949.7693
$ws.Range("K31").Value2 = 2949.7693
$ws.Range("M31").Value2 = -2654.7693

$ws.Range("H32").Value2 = 3083.3333
$ws.Range("I32").Value2 = 2000
$ws.Range("J32").Value2 = 3625
$ws.Range("K32").Value2 = 2000
$ws.Range("L32").Value2 = 3625
$ws.Range("M32").Value2 = -1684
$ws.Range("N32").Value2 = -4257

$ws.Range("H34").Value2 = 8250.117
$ws.Range("I34").Value2 = 2949.7693
$ws.Range("K34").Value2 = 2949.7693
$ws.Range("M34").Value2 = -2747.7693

$ws.Range("H35").Value2 = 340.65216
$ws.Range("I35").Value2 = 150.45
$ws.Range("J35").Value2 = 1608.6666
$ws.Range("K35").Value2 = 150.45
$ws.Range("L35").Value2 = 1608.6666
$ws.Range("M35").Value2 = 143.55
$ws.Range("N35").Value2 = -2196.6666

$ws.Range("H36").Value2 = 36163.332
$ws.Range("I36").Value2 = 0
$ws.Range("J36").Value2 = 36163.332
$ws.Range("K36").Value2 = 0
$ws.Range("L36").Value2 = 36163.332
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value2 = -36939.332

$ws.Range("H37").Value2 = 6700
$ws.Range("J37").Value2 = 6700
$ws.Range("L37").Value2 = 6700
$ws.Range("N37").Value2 = -6914

$ws.Range("H40").Value2 = 36163.332
$ws.Range("I40").Value2 = 0
$ws.Range("J40").Value2 = 36163.332
$ws.Range("K40").Value2 = 0
$ws.Range("L40").Value2 = 36163.332
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value2 = -36483.332

$ws.Range("H134").Value2 = 6213.4688
$ws.Range("I134").Value2 = 2414.077
$ws.Range("K134").Value2 = 7242.231000000001
$ws.Range("M134").Value2 = -4707.231000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value2 = 33333584
$ws.Range("I33").Value2 = 66666730
$ws.Range("J33").Value2 = 440.2
$ws.Range("K33").Value2 = 400000380
$ws.Range("L33").Value2 = 2641.2
$ws.Range("M33").Value2 = -400000097
$ws.Range("N33").Value2 = -3207.2

$ws.Range("H82").Value2 = 39250
$ws.Range("I82").Value2 = 20000
$ws.Range("J82").Value2 = 58500
$ws.Range("K82").Value2 = 60000
$ws.Range("L82").Value2 = 175500
$ws.Range("M82").Value2 = -59594
$ws.Range("N82").Value2 = -176312

$ws.Range("H85").Value2 = 39250
$ws.Range("I85").Value2 = 20000
$ws.Range("J85").Value2 = 58500
$ws.Range("K85").Value2 = 60000
$ws.Range("L85").Value2 = 175500
$ws.Range("M85").Value2 = -58596
$ws.Range("N85").Value2 = -178308

$ws.Range("H130").Value2 = 2166
$ws.Range("I130").Value2 = 2166
$ws.Range("K130").Value2 = 6498
$ws.Range("M130").Value2 = -1478

$ws.Range("H137").Value2 = 201742.2
$ws.Range("I137").Value2 = 144520.86
$ws.Range("K137").Value2 = 433562.58
$ws.Range("M137").Value2 = -428462.58

$ws.Range("H139").Value2 = 56979
$ws.Range("I139").Value2 = 65785.5
$ws.Range("K139").Value2 = 197356.5
$ws.Range("M139").Value2 = -192216.5

$ws.Range("H141").Value2 = 5610.5
$ws.Range("I141").Value2 = 5610.5
$ws.Range("K141").Value2 = 16831.5
$ws.Range("M141").Value2 = -11651.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 355854.66
$ws.Range("I70").Value2 = 731515.0600000001
$ws.Range("J70").Value2 = 11499.25
$ws.Range("K70").Value2 = 731515.0600000001
$ws.Range("L70").Value2 = 11499.25
$ws.Range("M70").Value2 = -731245.0600000001
$ws.Range("N70").Value2 = -12039.25

$ws.Range("H73").Value2 = 355854.66
$ws.Range("I73").Value2 = 731515.0600000001
$ws.Range("J73").Value2 = 11499.25
$ws.Range("K73").Value2 = 731515.0600000001
$ws.Range("L73").Value2 = 11499.25
$ws.Range("M73").Value2 = -730579.0600000001
$ws.Range("N73").Value2 = -13371.25

$ws.Range("H80").Value2 = 202846.4
$ws.Range("I80").Value2 = 3113
$ws.Range("J80").Value2 = 336002
$ws.Range("K80").Value2 = 3113
$ws.Range("L80").Value2 = 336002
$ws.Range("M80").Value2 = -2115
$ws.Range("N80").Value2 = -337998

$ws.Range("H83").Value2 = 202846.4
$ws.Range("I83").Value2 = 3113
$ws.Range("J83").Value2 = 336002
$ws.Range("K83").Value2 = 15565
$ws.Range("L83").Value2 = 1680010
$ws.Range("M83").Value2 = -10573
$ws.Range("N83").Value2 = -1689994

$ws.Range("H102").Value2 = 2995.8696
$ws.Range("I102").Value2 = 2882.75
$ws.Range("K102").Value2 = 2882.75
$ws.Range("M102").Value2 = -1260.75

$ws.Range("H122").Value2 = 1959994.2
$ws.Range("I122").Value2 = 2588814.5
$ws.Range("J122").Value2 = 3665
$ws.Range("K122").Value2 = 7766443.5
$ws.Range("L122").Value2 = 10995
$ws.Range("M122").Value2 = -7763993.5
$ws.Range("N122").Value2 = -15895

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 4301.8486
$ws.Range("I7").Value2 = 2359.0557
$ws.Range("J7").Value2 = 6633.2
$ws.Range("K7").Value2 = 2359.0557
$ws.Range("L7").Value2 = 6633.2
$ws.Range("M7").Value2 = -2247.0557
$ws.Range("N7").Value2 = -6857.2

$ws.Range("H122").Value2 = 3994.087
$ws.Range("I122").Value2 = 2704.6875
$ws.Range("J122").Value2 = 6941.2856
$ws.Range("K122").Value2 = 8114.0625
$ws.Range("L122").Value2 = 20823.8568
$ws.Range("M122").Value2 = -5664.0625
$ws.Range("N122").Value2 = -25723.8568

$ws.Range("H126").Value2 = 4301.8486
$ws.Range("I126").Value2 = 2359.0557
$ws.Range("J126").Value2 = 6633.2
$ws.Range("K126").Value2 = 7077.1671
$ws.Range("L126").Value2 = 19899.6
$ws.Range("M126").Value2 = -4607.1671
$ws.Range("N126").Value2 = -24839.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value2 = 10000
$ws.Range("J26").Value2 = 10000
$ws.Range("L26").Value2 = 10000
$ws.Range("N26").Value2 = -10586

$ws.Range("H43").Value2 = 8000
$ws.Range("I43").Value2 = 8000
$ws.Range("K43").Value2 = 8000
$ws.Range("M43").Value2 = -7851
